# Split the python-dict-like text that was stuffed into column B into
# separate, properly labeled columns (firstname/lastname/address1..4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------
# B1 already carries the bold/border header style (style index 1).
# Copy that formatting across the new header cells C1:G1, then set text.
$ws.Range("B1").Copy()
$ws.Range("C1:G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "firstname"
$ws.Range("C1").Value = "lastname"
$ws.Range("D1").Value = "address1"
$ws.Range("E1").Value = "address2"
$ws.Range("F1").Value = "address3"
$ws.Range("G1").Value = "address4"

# ---- Helper to write a text value without Excel mangling it --------
# Values like "09119" or " 08075" would otherwise be auto-converted to
# numbers (losing leading zeros / spaces). Forcing the number format to
# Text ("@") before the assignment keeps the literal string, and
# ClearFormats afterwards removes the temporary formatting again so the
# cell ends up with no explicit style, just like the rest of the data.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---- Data rows -------------------------------------------------------
$rows = @(
    @{r=2; firstname='John';                   lastname='Doe';      address1='120 jefferson st.';                address2='Riverside';   address3=' NJ'; address4=' 08075'},
    @{r=3; firstname='Jack';                   lastname='McGinnis'; address1='220 hobo Av.';                      address2='Phila';       address3=' PA'; address4='09119'},
    @{r=4; firstname='John "Da Man"';          lastname='Repici';   address1='120 Jefferson St.';                 address2='Riverside';   address3=' NJ'; address4='08075'},
    @{r=5; firstname='Stephen';                lastname='Tyler';    address1='7452 Terrace "At the Plaza" road';  address2='SomeTown';    address3='SD';  address4=' 91234'},
    @{r=6; firstname='';                       lastname='Blankman'; address1='';                                  address2='SomeTown';    address3=' SD'; address4=' 00298'},
    @{r=7; firstname='Joan "the bone", Anne';  lastname='Jet';      address1='9th, at Terrace plc';               address2='Desert City'; address3='CO';  address4='00123'}
)

foreach ($row in $rows) {
    $r = $row.r

    Set-TextValue $ws.Cells.Item($r, 2) $row.firstname
    Set-TextValue $ws.Cells.Item($r, 3) $row.lastname
    Set-TextValue $ws.Cells.Item($r, 4) $row.address1
    Set-TextValue $ws.Cells.Item($r, 5) $row.address2
    Set-TextValue $ws.Cells.Item($r, 6) $row.address3
    Set-TextValue $ws.Cells.Item($r, 7) $row.address4
}
